$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1335.6666
$ws.Range("J19").Value = 2097.5
$ws.Range("L19").Value = 2097.5
$ws.Range("N19").Value = -2447.5
$ws.Range("H40").Value = 1045.871
$ws.Range("I40").Value = 866.26086
$ws.Range("K40").Value = 866.26086
$ws.Range("M40").Value = -691.26086
$ws.Range("H62").Value = 10000
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 10000
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H98").Value = 1691.0714
$ws.Range("I98").Value = 1759.5834
$ws.Range("K98").Value = 1759.5834
$ws.Range("M98").Value = -261.5834
$ws.Range("H122").Value = 1691.0714
$ws.Range("I122").Value = 1759.5834
$ws.Range("K122").Value = 5278.7502
$ws.Range("M122").Value = -2828.7502
$ws.Range("H135").Value = 206.2
$ws.Range("I135").Value = 77.333336
$ws.Range("K135").Value = 696.0000240000001
$ws.Range("M135").Value = 1838.999976
$ws.Range("H138").Value = 2185.61
$ws.Range("I138").Value = 986.3333
$ws.Range("J138").Value = 2304.2197
$ws.Range("K138").Value = 2958.9999
$ws.Range("L138").Value = 6912.659100000001
$ws.Range("M138").Value = 2181.0001
$ws.Range("N138").Value = -17192.6591

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1369.1765
$ws.Range("I61").Value = 1183
$ws.Range("J61").Value = 2238
$ws.Range("K61").Value = 1183
$ws.Range("L61").Value = 2238
$ws.Range("M61").Value = -971
$ws.Range("N61").Value = -2662
$ws.Range("H74").Value = 1124
$ws.Range("I74").Value = 1136
$ws.Range("J74").Value = 1076
$ws.Range("K74").Value = 1136
$ws.Range("L74").Value = 1076
$ws.Range("M74").Value = -262
$ws.Range("N74").Value = -2824
$ws.Range("H77").Value = 1124
$ws.Range("I77").Value = 1136
$ws.Range("K77").Value = 5680
$ws.Range("L77").Value = 5380
$ws.Range("M77").Value = -1312
$ws.Range("N77").Value = -14116
$ws.Range("H102").Value = 41692620
$ws.Range("I102").Value = 83383340
$ws.Range("K102").Value = 83383340
$ws.Range("M102").Value = -83381718
$ws.Range("H103").Value = 63000
$ws.Range("J103").Value = 63000
$ws.Range("L103").Value = 63000
$ws.Range("N103").Value = -65344
$ws.Range("H132").Value = 2600.6765
$ws.Range("I132").Value = 2353.1853
$ws.Range("J132").Value = 3555.2856
$ws.Range("K132").Value = 7059.5559
$ws.Range("L132").Value = 10665.8568
$ws.Range("M132").Value = -4529.5559
$ws.Range("N132").Value = -15725.8568
$ws.Range("H136").Value = 1369.1765
$ws.Range("I136").Value = 1183
$ws.Range("J136").Value = 2238
$ws.Range("K136").Value = 3549
$ws.Range("L136").Value = 6714
$ws.Range("M136").Value = -999
$ws.Range("N136").Value = -11814

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 62500524
$ws.Range("I94").Value = 62500524
$ws.Range("K94").Value = 62500524
$ws.Range("M94").Value = -62500073
$ws.Range("H134").Value = 1697.25
$ws.Range("I134").Value = 1674.0667
$ws.Range("J134").Value = 1766.8
$ws.Range("K134").Value = 5022.2001
$ws.Range("L134").Value = 5300.4
$ws.Range("M134").Value = -2487.2001
$ws.Range("N134").Value = -10370.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4257797
$ws.Range("I62").Value = 2476.7441
$ws.Range("J62").Value = 50002490
$ws.Range("K62").Value = 2476.7441
$ws.Range("L62").Value = 50002490
$ws.Range("M62").Value = -1852.7441
$ws.Range("N62").Value = -50003738
$ws.Range("H65").Value = 4257797
$ws.Range("I65").Value = 2476.7441
$ws.Range("J65").Value = 50002490
$ws.Range("K65").Value = 12383.7205
$ws.Range("L65").Value = 250012450
$ws.Range("M65").Value = -9263.720499999999
$ws.Range("N65").Value = -250018690
$ws.Range("H99").Value = 1895.6666
$ws.Range("I99").Value = 1915.1111
$ws.Range("J99").Value = 1866.5
$ws.Range("K99").Value = 1915.1111
$ws.Range("L99").Value = 1866.5
$ws.Range("M99").Value = -417.1111000000001
$ws.Range("N99").Value = -4862.5
$ws.Range("H110").Value = 38000
$ws.Range("J110").Value = 38000
$ws.Range("L110").Value = 38000
$ws.Range("N110").Value = -46180
$ws.Range("H126").Value = 1895.6666
$ws.Range("I126").Value = 1915.1111
$ws.Range("J126").Value = 1866.5
$ws.Range("K126").Value = 5745.3333
$ws.Range("L126").Value = 5599.5
$ws.Range("M126").Value = -3275.3333
$ws.Range("N126").Value = -10539.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1579.2142
$ws.Range("I5").Value = 1460.7
$ws.Range("J5").Value = 1875.5
$ws.Range("K5").Value = 4382.1
$ws.Range("L5").Value = 5626.5
$ws.Range("M5").Value = -4270.1
$ws.Range("N5").Value = -5850.5
$ws.Range("H56").Value = 3929.75
$ws.Range("I56").Value = 3929.75
$ws.Range("K56").Value = 3929.75
$ws.Range("M56").Value = -3399.75
$ws.Range("H75").Value = 2186
$ws.Range("I75").Value = 500
$ws.Range("J75").Value = 2607.5
$ws.Range("K75").Value = 1500
$ws.Range("L75").Value = 7822.5
$ws.Range("M75").Value = -502
$ws.Range("N75").Value = -9818.5
$ws.Range("H78").Value = 2186
$ws.Range("I78").Value = 500
$ws.Range("J78").Value = 2607.5
$ws.Range("K78").Value = 4500
$ws.Range("L78").Value = 23467.5
$ws.Range("M78").Value = 492
$ws.Range("N78").Value = -33451.5
$ws.Range("H81").Value = 1325
$ws.Range("I81").Value = 1325
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 3975
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -2852
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 1325
$ws.Range("I84").Value = 1325
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 11925
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -6309
$ws.Range("N84").ClearContents()
$ws.Range("H135").Value = 1579.2142
$ws.Range("I135").Value = 1460.7
$ws.Range("J135").Value = 1875.5
$ws.Range("K135").Value = 13146.3
$ws.Range("L135").Value = 16879.5
$ws.Range("M135").Value = -10611.3
$ws.Range("N135").Value = -21949.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1000
$ws.Range("I97").Value = 1000
$ws.Range("J97").Value = 1000
$ws.Range("K97").Value = 1000
$ws.Range("L97").Value = 1000
$ws.Range("M97").Value = -504
$ws.Range("N97").Value = -1992
$ws.Range("H126").Value = 2707.7273
$ws.Range("I126").Value = 1797
$ws.Range("K126").Value = 5391
$ws.Range("M126").Value = -2921
$ws.Range("H132").Value = 1951.6
$ws.Range("I132").Value = 1551.2142
$ws.Range("J132").Value = 3553.1428
$ws.Range("K132").Value = 4653.642599999999
$ws.Range("L132").Value = 10659.4284
$ws.Range("M132").Value = -2123.642599999999
$ws.Range("N132").Value = -15719.4284

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1319.5625
$ws.Range("I68").Value = 1173.2174
$ws.Range("J68").Value = 1693.5555
$ws.Range("K68").Value = 1173.2174
$ws.Range("L68").Value = 1693.5555
$ws.Range("M68").Value = -424.2174
$ws.Range("N68").Value = -3191.5555
$ws.Range("H71").Value = 1319.5625
$ws.Range("I71").Value = 1173.2174
$ws.Range("J71").Value = 1693.5555
$ws.Range("K71").Value = 5866.087
$ws.Range("L71").Value = 8467.7775
$ws.Range("M71").Value = -2122.087
$ws.Range("N71").Value = -15955.7775
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
$ws.Range("H132").Value = 47904.5
$ws.Range("I132").Value = 2229.8
$ws.Range("K132").Value = 6689.400000000001
$ws.Range("M132").Value = -4159.400000000001
$ws.Range("H136").Value = 1831.25
$ws.Range("I136").Value = 1358.3334
$ws.Range("J136").Value = 3250
$ws.Range("K136").Value = 4075.0002
$ws.Range("L136").Value = 9750
$ws.Range("M136").Value = -1525.0002
$ws.Range("N136").Value = -14850

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 9390
$ws.Range("J92").Value = 9390
$ws.Range("L92").Value = 9390
$ws.Range("N92").Value = -14382
$ws.Range("H100").Value = 475.73685
$ws.Range("I100").Value = 495
$ws.Range("J100").Value = 442.7143
$ws.Range("K100").Value = 990
$ws.Range("L100").Value = 885.4286
$ws.Range("M100").Value = -449
$ws.Range("N100").Value = -1967.4286
$ws.Range("H126").Value = 83334776
$ws.Range("I126").Value = 125000850
$ws.Range("K126").Value = 375002550
$ws.Range("M126").Value = -375000080
$ws.Range("H136").Value = 1352.2727
$ws.Range("I136").Value = 1016.06665
$ws.Range("K136").Value = 3048.19995
$ws.Range("M136").Value = -498.1999500000002

Write-Output "Applied 232 cell changes"